$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "29.978.14"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "  +0.52%  "
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "1.910.14"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "  +0.93%  "
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "0.9999"
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "  +0.06%  "
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "0.8046"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "  +5.07%  "
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "242.05"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "  +1.18%  "
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "  -0.12%  "
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.3164"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "  +3.92%  "
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "26.49"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "  +4.52%  "
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.06914"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.07990"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "  +0.00%  "
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "1.914.65"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "  +1.28%  "
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "0.7408"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "  -1.03%  "
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "5.203"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "  +0.10%  "
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "93.06"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "  +2.26%  "
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "29.981.50"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "  +0.51%  "
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "14.02"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "  +1.07%  "
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "5.888"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "  -1.14%  "
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "246.53"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = "  +5.11%  "
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "0.000007743"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = "  +1.02%  "
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.9998"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = "  +0.01%  "
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "2.152.70"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "  +0.89%  "
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "0.9997"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = "  +0.05%  "
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "6.849"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = "  -1.11%  "
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "167.82"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = "  +1.48%  "
$ws.Cells.Item(25,5).Style = "Normal"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "9.225"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = "  -0.17%  "
$ws.Cells.Item(26,5).Style = "Normal"
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "0.1422"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value = "  +9.23%  "
$ws.Cells.Item(27,5).Style = "Normal"
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "18.93"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value = "  +1.25%  "
$ws.Cells.Item(28,5).Style = "Normal"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "2.039"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = "  -0.16%  "
$ws.Cells.Item(29,5).Style = "Normal"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "1.365"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value = "  +1.70%  "
$ws.Cells.Item(30,5).Style = "Normal"
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "1.514"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,5).Value = "  +0.19%  "
$ws.Cells.Item(31,5).Style = "Normal"
$ws.Cells.Item(32,5).NumberFormat = "@"
$ws.Cells.Item(32,5).Value = "  +0.97%  "
$ws.Cells.Item(32,5).Style = "Normal"
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "4.086"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).NumberFormat = "@"
$ws.Cells.Item(33,5).Value = "  +1.59%  "
$ws.Cells.Item(33,5).Style = "Normal"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "0.05479"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).NumberFormat = "@"
$ws.Cells.Item(34,5).Value = "  +2.16%  "
$ws.Cells.Item(34,5).Style = "Normal"
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "1.265"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).NumberFormat = "@"
$ws.Cells.Item(35,5).Value = "  +1.38%  "
$ws.Cells.Item(35,5).Style = "Normal"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.7336"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).NumberFormat = "@"
$ws.Cells.Item(36,5).Value = "  +0.99%  "
$ws.Cells.Item(36,5).Style = "Normal"
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "2.712"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).NumberFormat = "@"
$ws.Cells.Item(37,5).Value = "  +0.03%  "
$ws.Cells.Item(37,5).Style = "Normal"
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "0.01924"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = "  +0.05%  "
$ws.Cells.Item(38,5).Style = "Normal"
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "2.789"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = "  +0.70%  "
$ws.Cells.Item(39,5).Style = "Normal"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "6.160"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = "  -0.33%  "
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.4422"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = "  +0.44%  "
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "72.49"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = "  +0.47%  "
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.9981"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = "  -0.18%  "
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.8359"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = "  +1.49%  "
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "1.877"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = "  -1.87%  "
$ws.Cells.Item(45,5).Style = "Normal"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "7.556"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = "  -0.32%  "
$ws.Cells.Item(46,5).Style = "Normal"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "100.43"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = "  -0.59%  "
$ws.Cells.Item(47,5).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = "  -0.49%  "
$ws.Cells.Item(48,5).Style = "Normal"
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "982.84"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,5).Value = "  +6.52%  "
$ws.Cells.Item(49,5).Style = "Normal"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "2.058.08"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,5).Value = "  +1.09%  "
$ws.Cells.Item(50,5).Style = "Normal"
$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,5).Value = "  +0.37%  "
$ws.Cells.Item(51,5).Style = "Normal"
